$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the empty bullet paragraph (ListParagraph, ilvl=1) that sits
#    between "...insert ID # here." and "Compare and contrast the
#    following:" so that the "Compare and contrast" bullet becomes the
#    immediate successor of the hyperlink bullet.
# ---------------------------------------------------------------------
for ($j = 1; $j -le $d.Paragraphs.Count; $j++) {
    $pp = $d.Paragraphs.Item($j)
    if ($pp.Range.Text.Trim() -eq "" -and $pp.Range.Text.Length -le 2) {
        $prev = $d.Paragraphs.Item($j - 1)
        $next = $d.Paragraphs.Item($j + 1)
        if ($prev.Range.Text.StartsWith("https://calvincs262-monopoly.appspot.com/monopoly/v1/player/") -and `
            $next.Range.Text.StartsWith("Compare and contrast the following")) {
            $pp.Range.Delete()
            break
        }
    }
}

# ---------------------------------------------------------------------
# 2. Locate "JDBC statement objects seem to be a wrapper for the SQL
#    query statements." and turn it into:
#      "JDBC statement objects seem to be a wrapper"
#      "/interface/adapter"
#      <bookmarkStart/End name="_GoBack"/>
#      " for the SQL query statements."
# ---------------------------------------------------------------------
$targetIdx = -1
for ($j = 1; $j -le $d.Paragraphs.Count; $j++) {
    $pp = $d.Paragraphs.Item($j)
    if ($pp.Range.Text.StartsWith("JDBC statement objects seem to be a wrapper for the SQL")) {
        $targetIdx = $j
        break
    }
}

$pp = $d.Paragraphs.Item($targetIdx)
$pStart = $pp.Range.Start
$full = $pp.Range.Text
$wrapperEnd = $pStart + $full.IndexOf("wrapper") + "wrapper".Length

# Insert "/interface/adapter" right after "wrapper"
$insertPoint = $d.Range($wrapperEnd, $wrapperEnd)
$insertPoint.InsertAfter("/interface/adapter")

$adapterStart = $wrapperEnd
$adapterEnd = $wrapperEnd + "/interface/adapter".Length

# Force a run boundary between "wrapper" and "/interface/adapter" (no
# visible formatting change — use a transient bookmark purely to split
# the backing run, then remove the transient bookmark again).
$tmpRange = $d.Range($adapterStart, $adapterStart)
$d.Bookmarks.Add("TmpSplitBoundary", $tmpRange)
$d.Bookmarks.Item("TmpSplitBoundary").Delete()

# ---------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark from its old position (inside the
#    "JDBC result sets..." paragraph) to right after "/interface/adapter"
#    in the paragraph edited above.
# ---------------------------------------------------------------------
$oldBm = $d.Bookmarks.Item("_GoBack")
$oldPos = $oldBm.Start
$oldBm.Delete()

# The old bookmark previously split "...Contains " | "the rows...";
# removing it should leave those two runs re-merged into a single run.
$boundary = $d.Range($oldPos - 1, $oldPos)
$ch = $boundary.Text
$boundary.Text = ""
$reinsert = $d.Range($oldPos - 1, $oldPos - 1)
$reinsert.InsertAfter($ch)

# Re-create "_GoBack" at the new location.
$newBmRange = $d.Range($adapterEnd, $adapterEnd)
$d.Bookmarks.Add("_GoBack", $newBmRange)
